$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$lightBlue = 15128749   # BGR-encoded 0xADD8E6 (fill color used for "Minor Deficiencies")

# Row 4 (1-based, data row "asdsa/asd/ads/Accepted/asd") - shading + text only
$cell4 = $table.Cell(4, 4)
$cell4.Shading.BackgroundPatternColor = $lightBlue
$cell4.Range.Text = "Minor Deficiencies"

# Row 5 (1-based, data row "new/king/das/Accepted/ads") - row height + text + shading + text
$row5 = $table.Rows.Item(5)
$row5.Height = 28.8

$cell5b = $table.Cell(5, 2)
$cell5b.Range.Text = "da"

$cell5 = $table.Cell(5, 4)
$cell5.Shading.BackgroundPatternColor = $lightBlue
$cell5.Range.Text = "Minor Deficiencies"
